$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.16"
$ws.Range("F2").Value = '28-12-2022'
$ws.Range("G2").Value = "'0"
$ws.Range("D3").Value = "'23.90"
$ws.Range("F3").Value = '28-12-2022'
$ws.Range("G3").Value = "'0"
$ws.Range("D4").Value = "'5.364"
$ws.Range("F4").Value = '28-12-2022'
$ws.Range("G4").Value = "'0"
$ws.Range("D5").Value = "'0.05809"
$ws.Range("F5").Value = '28-12-2022'
$ws.Range("G5").Value = "'0"
$ws.Range("D6").Value = "'3.370"
$ws.Range("F6").Value = '28-12-2022'
$ws.Range("G6").Value = "'0"
$ws.Range("D7").Value = "'6.479"
$ws.Range("F7").Value = '28-12-2022'
$ws.Range("G7").Value = "'0"
$ws.Range("D8").Value = "'0.8107"
$ws.Range("F8").Value = '28-12-2022'
$ws.Range("G8").Value = "'0"
$ws.Range("F9").Value = '28-12-2022'
$ws.Range("G9").Value = "'0"
$ws.Range("D10").Value = "'0.1399"
$ws.Range("F10").Value = '28-12-2022'
$ws.Range("G10").Value = "'0"
$ws.Range("D11").Value = "'0.07387"
$ws.Range("F11").Value = '28-12-2022'
$ws.Range("G11").Value = "'0"
$ws.Range("D12").Value = "'0.03210"
$ws.Range("F12").Value = '28-12-2022'
$ws.Range("G12").Value = "'0"
$ws.Range("D13").Value = "'0.03076"
$ws.Range("F13").Value = '28-12-2022'
$ws.Range("G13").Value = "'0"
$ws.Range("D14").Value = "'0.09362"
$ws.Range("F14").Value = '28-12-2022'
$ws.Range("G14").Value = "'0"
$ws.Range("D15").Value = "'3.847"
$ws.Range("F15").Value = '28-12-2022'
$ws.Range("G15").Value = "'0"
$ws.Range("D16").Value = "'0.001547"
$ws.Range("F16").Value = '28-12-2022'
$ws.Range("G16").Value = "'0"
$ws.Range("D17").Value = "'0.04697"
$ws.Range("F17").Value = '28-12-2022'
$ws.Range("G17").Value = "'0"
$ws.Range("D18").Value = "'0.0005977"
$ws.Range("F18").Value = '28-12-2022'
$ws.Range("G18").Value = "'0"
$ws.Range("D19").Value = "'0.005922"
$ws.Range("F19").Value = '28-12-2022'
$ws.Range("G19").Value = "'0"
$ws.Range("D20").Value = "'0.001255"
$ws.Range("F20").Value = '28-12-2022'
$ws.Range("G20").Value = "'0"
$ws.Range("D21").Value = "'0.004670"
$ws.Range("F21").Value = '28-12-2022'
$ws.Range("G21").Value = "'0"
$ws.Range("D22").Value = "'0.00008795"
$ws.Range("E22").Value = '21NitroExNTXBestin24h'
$ws.Range("F22").Value = '28-12-2022'
$ws.Range("G22").Value = "'0"
$ws.Range("D23").Value = "'3.594"
$ws.Range("F23").Value = '28-12-2022'
$ws.Range("G23").Value = "'0"
$ws.Range("F24").Value = '28-12-2022'
$ws.Range("G24").Value = "'0"
$ws.Range("D25").Value = "'0.3176"
$ws.Range("F25").Value = '28-12-2022'
$ws.Range("G25").Value = "'0"
$ws.Range("D26").Value = "'0.1318"
$ws.Range("F26").Value = '28-12-2022'
$ws.Range("G26").Value = "'0"
$ws.Range("F27").Value = '28-12-2022'
$ws.Range("G27").Value = "'0"
$ws.Range("F28").Value = '28-12-2022'
$ws.Range("G28").Value = "'0"
$ws.Range("F29").Value = '28-12-2022'
$ws.Range("G29").Value = "'0"
$ws.Range("F30").Value = '28-12-2022'
$ws.Range("G30").Value = "'0"
$ws.Range("F31").Value = '28-12-2022'
$ws.Range("G31").Value = "'0"
$ws.Range("F32").Value = '28-12-2022'
$ws.Range("G32").Value = "'0"
$ws.Range("F33").Value = '28-12-2022'
$ws.Range("G33").Value = "'0"
$ws.Range("F34").Value = '28-12-2022'
$ws.Range("G34").Value = "'0"
$ws.Range("F35").Value = '28-12-2022'
$ws.Range("G35").Value = "'0"
$ws.Range("F36").Value = '28-12-2022'
$ws.Range("G36").Value = "'0"
$ws.Range("F37").Value = '28-12-2022'
$ws.Range("G37").Value = "'0"
$ws.Range("F38").Value = '28-12-2022'
$ws.Range("G38").Value = "'0"
$ws.Range("F39").Value = '28-12-2022'
$ws.Range("G39").Value = "'0"
$ws.Range("D40").Value = "'0.03853"
$ws.Range("F40").Value = '28-12-2022'
$ws.Range("G40").Value = "'0"
$ws.Range("B41").Value = 'CEJI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D41").Value = "'0.002758"
$ws.Range("E41").Value = '40CEJICEJI'
$ws.Range("F41").Value = '28-12-2022'
$ws.Range("G41").Value = "'0"
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").Value = "'0.003055"
$ws.Range("E42").Value = '41KickTokenKICKWorstin24h'
$ws.Range("F42").Value = '28-12-2022'
$ws.Range("G42").Value = "'0"
$ws.Range("D43").Value = "'0.1065"
$ws.Range("F43").Value = '28-12-2022'
$ws.Range("G43").Value = "'0"
$ws.Range("D44").Value = "'0.009062"
$ws.Range("F44").Value = '28-12-2022'
$ws.Range("G44").Value = "'0"
$ws.Range("F45").Value = '28-12-2022'
$ws.Range("G45").Value = "'0"
$ws.Range("F46").Value = '28-12-2022'
$ws.Range("G46").Value = "'0"
$ws.Range("F47").Value = '28-12-2022'
$ws.Range("G47").Value = "'0"
$ws.Range("D48").Value = "'0.001847"
$ws.Range("E48").Value = '47BOLOBOLO'
$ws.Range("F48").Value = '28-12-2022'
$ws.Range("G48").Value = "'0"
$ws.Range("F49").Value = '28-12-2022'
$ws.Range("G49").Value = "'0"
$ws.Range("F50").Value = '28-12-2022'
$ws.Range("G50").Value = "'0"
$ws.Range("F51").Value = '28-12-2022'
$ws.Range("G51").Value = "'0"
